$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# Row 64 used to be a lone, orphaned tag cell (C64 only). The commit
# turns it into a full row for LeetCode #56 "Merge Intervals", reusing
# the formatting of the row right above it (row 63).
# ------------------------------------------------------------------
$ws.Range("A63:I63").Copy() | Out-Null
$ws.Range("A64:I64").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A64").Value = 56
$ws.Range("B64").Value = "Merge Intervals"
$ws.Range("C64").Value = "#array #sorting #核心 "
$ws.Range("D64").Value = "medium"
$ws.Range("E64").Value = 4
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 45846
$ws.Range("I64").Value = 45846
$ws.Range("A64:I64").RowHeight = 34

# ------------------------------------------------------------------
# Row 65: brand-new entry for LeetCode #3439 "Reschedule Meetings for
# Maximum Free Time I" (the problem the new markdown note documents).
# Row 61 already has the full A:J formatting (including the "?" J
# column marker), so borrow that as the template.
# ------------------------------------------------------------------
$ws.Range("A61:J61").Copy() | Out-Null
$ws.Range("A65:J65").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A65").Value = 3439
$ws.Range("B65").Value = "Reschedule Meetings for Maximum Free Time I"
$ws.Range("C65").Value = "#array #greedy #sliding-window "
$ws.Range("D65").Value = "medium"
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 1
$ws.Range("G65").Value = 45
$ws.Range("H65").Value = 45847
$ws.Range("I65").Value = 45847
$ws.Range("J65").Value = "?"
$ws.Range("A65:J65").RowHeight = 34

# ------------------------------------------------------------------
# Selection / scroll bookkeeping, mirroring the saved view state.
# ------------------------------------------------------------------
$ws.Range("J65").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 60
